$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.925.47'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '''1.845.68'
$ws.Range("E3").Value = '  +1.31%  '
$ws.Range("D4").Value = '''1.005'
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").Value = '''309.22'
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("D7").Value = '''0.4722'
$ws.Range("E7").Value = '  +0.78%  '
$ws.Range("D8").Value = '''0.3681'
$ws.Range("E8").Value = '  +2.34%  '
$ws.Range("D9").Value = '''0.07212'
$ws.Range("E9").Value = '  +1.20%  '
$ws.Range("D10").Value = '''0.9252'
$ws.Range("E10").Value = '  +2.56%  '
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").Value = '''0.07622'
$ws.Range("E12").Value = '  -2.47%  '
$ws.Range("D13").Value = '''1.891.24'
$ws.Range("E13").Value = '  +4.44%  '
$ws.Range("D14").Value = '''5.314'
$ws.Range("E14").Value = '  +1.09%  '
$ws.Range("D15").Value = '''6.405'
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("D16").Value = '''88.58'
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").Value = '''0.000008668'
$ws.Range("E18").Value = '  +1.25%  '
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D20").Value = '''26.950.92'
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '''14.56'
$ws.Range("E21").Value = '  +2.65%  '
$ws.Range("D22").Value = '''5.036'
$ws.Range("E22").Value = '  +0.43%  '
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("D24").Value = '''1.917'
$ws.Range("E24").Value = '  -0.99%  '
$ws.Range("D25").Value = '''152.13'
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("D26").Value = '''18.15'
$ws.Range("E26").Value = '  +1.32%  '
$ws.Range("D27").Value = '''2.007'
$ws.Range("E27").Value = '  +1.71%  '
$ws.Range("D28").Value = '''114.25'
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("D29").Value = '''4.947'
$ws.Range("E29").Value = '  +2.92%  '
$ws.Range("D30").Value = '''0.08835'
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("D31").Value = '''3.310'
$ws.Range("E31").Value = '  +5.29%  '
$ws.Range("D32").Value = '''0.7469'
$ws.Range("E32").Value = '  +2.22%  '
$ws.Range("D33").Value = '''1.170'
$ws.Range("E33").Value = '  +4.17%  '
$ws.Range("D34").Value = '''2.776'
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").Value = '''4.488'
$ws.Range("E35").Value = '  +1.02%  '
$ws.Range("D36").Value = '''1.089'
$ws.Range("E36").Value = '  +1.10%  '
$ws.Range("D37").Value = '''0.05264'
$ws.Range("E37").Value = '  +2.89%  '
$ws.Range("D38").Value = '''0.01949'
$ws.Range("E38").Value = '  +1.10%  '
$ws.Range("D39").Value = '''2.963'
$ws.Range("E39").Value = '  +1.48%  '
$ws.Range("D40").Value = '''0.5228'
$ws.Range("E40").Value = '  +3.34%  '
$ws.Range("D41").Value = '''6.924'
$ws.Range("E41").Value = '  +1.50%  '
$ws.Range("D42").Value = '''0.1514'
$ws.Range("E42").Value = '  +1.13%  '
$ws.Range("D43").Value = '''8.229'
$ws.Range("E43").Value = '  +2.90%  '
$ws.Range("D44").Value = '''10.50'
$ws.Range("E44").Value = '  +4.91%  '
$ws.Range("D45").Value = '''0.4705'
$ws.Range("E45").Value = '  +0.61%  '
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("D47").Value = '''101.98'
$ws.Range("E47").Value = '  +2.98%  '
$ws.Range("D48").Value = '''1.608'
$ws.Range("E48").Value = '  +3.01%  '
$ws.Range("D49").Value = '''65.37'
$ws.Range("E49").Value = '  +2.48%  '
$ws.Range("D50").Value = '''0.06031'
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").Value = '''0.8862'
$ws.Range("E51").Value = '  +4.15%  '
